# Add ANR (Advanced Nuclear Reactor) electricity and hydrogen production
# plant-type rows to the NewTechFramework sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- ANR electricity-producing plant types (rows 18-22) ---

# Row 18: iPWR
$ws.Cells.Item(18, 1).Value = "iPWR"
$ws.Cells.Item(18, 2).Value = "ANRElec"
$ws.Cells.Item(18, 3).Value = "NA"
$ws.Cells.Item(18, 4).Value = "Nuclear Fuel"
$ws.Cells.Item(18, 5).Value = "thermal"
$ws.Cells.Item(18, 6).Value = 77
$ws.Cells.Item(18, 11).Value = "Yes"
$ws.Cells.Item(18, 12).Value = 0
$ws.Cells.Item(18, 13).Value = 0
$ws.Cells.Item(18, 14).Value = 0
$ws.Cells.Item(18, 15).Value = 60
$ws.Cells.Item(18, 16).Value = 0
$ws.Cells.Item(18, 17).Value = "NA"
$ws.Cells.Item(18, 18).Value = "NA"
$ws.Cells.Item(18, 19).Value = "NA"
$ws.Cells.Item(18, 20).Value = "NA"
$ws.Cells.Item(18, 21).Value = "NA"
$ws.Cells.Item(18, 22).Value = "NA"
$ws.Cells.Item(18, 23).Value = "NA"
$ws.Cells.Item(18, 24).Value = "Nuclear"

# Row 19: HTGR
$ws.Cells.Item(19, 1).Value = "HTGR"
$ws.Cells.Item(19, 2).Value = "ANRElec"
$ws.Cells.Item(19, 3).Value = "NA"
$ws.Cells.Item(19, 4).Value = "Nuclear Fuel"
$ws.Cells.Item(19, 5).Value = "thermal"
$ws.Cells.Item(19, 6).Value = 164
$ws.Cells.Item(19, 11).Value = "Yes"
$ws.Cells.Item(19, 12).Value = 0
$ws.Cells.Item(19, 13).Value = 0
$ws.Cells.Item(19, 14).Value = 0
$ws.Cells.Item(19, 15).Value = 60
$ws.Cells.Item(19, 16).Value = 0
$ws.Cells.Item(19, 17).Value = "NA"
$ws.Cells.Item(19, 18).Value = "NA"
$ws.Cells.Item(19, 19).Value = "NA"
$ws.Cells.Item(19, 20).Value = "NA"
$ws.Cells.Item(19, 21).Value = "NA"
$ws.Cells.Item(19, 22).Value = "NA"
$ws.Cells.Item(19, 23).Value = "NA"
$ws.Cells.Item(19, 24).Value = "Nuclear"

# Row 20: PBRHTGR
$ws.Cells.Item(20, 1).Value = "PBRHTGR"
$ws.Cells.Item(20, 2).Value = "ANRElec"
$ws.Cells.Item(20, 3).Value = "NA"
$ws.Cells.Item(20, 4).Value = "Nuclear Fuel"
$ws.Cells.Item(20, 5).Value = "thermal"
$ws.Cells.Item(20, 6).Value = 80
$ws.Cells.Item(20, 11).Value = "Yes"
$ws.Cells.Item(20, 12).Value = 0
$ws.Cells.Item(20, 13).Value = 0
$ws.Cells.Item(20, 14).Value = 0
$ws.Cells.Item(20, 15).Value = 60
$ws.Cells.Item(20, 16).Value = 0
$ws.Cells.Item(20, 17).Value = "NA"
$ws.Cells.Item(20, 18).Value = "NA"
$ws.Cells.Item(20, 19).Value = "NA"
$ws.Cells.Item(20, 20).Value = "NA"
$ws.Cells.Item(20, 21).Value = "NA"
$ws.Cells.Item(20, 22).Value = "NA"
$ws.Cells.Item(20, 23).Value = "NA"
$ws.Cells.Item(20, 24).Value = "Nuclear"

# Row 21: iMSR
$ws.Cells.Item(21, 1).Value = "iMSR"
$ws.Cells.Item(21, 2).Value = "ANRElec"
$ws.Cells.Item(21, 3).Value = "NA"
$ws.Cells.Item(21, 4).Value = "Nuclear Fuel"
$ws.Cells.Item(21, 5).Value = "thermal"
$ws.Cells.Item(21, 6).Value = 141
$ws.Cells.Item(21, 11).Value = "Yes"
$ws.Cells.Item(21, 12).Value = 0
$ws.Cells.Item(21, 13).Value = 0
$ws.Cells.Item(21, 14).Value = 0
$ws.Cells.Item(21, 15).Value = 60
$ws.Cells.Item(21, 16).Value = 0
$ws.Cells.Item(21, 17).Value = "NA"
$ws.Cells.Item(21, 18).Value = "NA"
$ws.Cells.Item(21, 19).Value = "NA"
$ws.Cells.Item(21, 20).Value = "NA"
$ws.Cells.Item(21, 21).Value = "NA"
$ws.Cells.Item(21, 22).Value = "NA"
$ws.Cells.Item(21, 23).Value = "NA"
$ws.Cells.Item(21, 24).Value = "Nuclear"

# Row 22: Micro
$ws.Cells.Item(22, 1).Value = "Micro"
$ws.Cells.Item(22, 2).Value = "ANRElec"
$ws.Cells.Item(22, 3).Value = "NA"
$ws.Cells.Item(22, 4).Value = "Nuclear Fuel"
$ws.Cells.Item(22, 5).Value = "thermal"
$ws.Cells.Item(22, 6).Value = 6.7
$ws.Cells.Item(22, 11).Value = "Yes"
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 13).Value = 0
$ws.Cells.Item(22, 14).Value = 0
$ws.Cells.Item(22, 15).Value = 60
$ws.Cells.Item(22, 16).Value = 0
$ws.Cells.Item(22, 17).Value = "NA"
$ws.Cells.Item(22, 18).Value = "NA"
$ws.Cells.Item(22, 19).Value = "NA"
$ws.Cells.Item(22, 20).Value = "NA"
$ws.Cells.Item(22, 21).Value = "NA"
$ws.Cells.Item(22, 22).Value = "NA"
$ws.Cells.Item(22, 23).Value = "NA"
$ws.Cells.Item(22, 24).Value = "Nuclear"

# --- ANR hydrogen (HTSE) production plant types (rows 23-27) ---

# Row 23: iPWRHTSE
$ws.Cells.Item(23, 1).Value = "iPWRHTSE"
$ws.Cells.Item(23, 2).Value = "ANRH2"
$ws.Cells.Item(23, 3).Value = "NA"
$ws.Cells.Item(23, 4).Value = "Nuclear Fuel"
$ws.Cells.Item(23, 5).Value = "h2"
$ws.Cells.Item(23, 6).Value = 77
$ws.Cells.Item(23, 7).Value = 100
$ws.Cells.Item(23, 11).Value = "Yes"
$ws.Cells.Item(23, 24).Value = "ANRH2"

# Row 24: HTGRHTSE
$ws.Cells.Item(24, 1).Value = "HTGRHTSE"
$ws.Cells.Item(24, 2).Value = "ANRH2"
$ws.Cells.Item(24, 3).Value = "NA"
$ws.Cells.Item(24, 4).Value = "Nuclear Fuel"
$ws.Cells.Item(24, 5).Value = "h2"
$ws.Cells.Item(24, 6).Value = 164
$ws.Cells.Item(24, 7).Value = 30
$ws.Cells.Item(24, 11).Value = "Yes"
$ws.Cells.Item(24, 24).Value = "ANRH2"

# Row 25: PBRHTGRHTSE
$ws.Cells.Item(25, 1).Value = "PBRHTGRHTSE"
$ws.Cells.Item(25, 2).Value = "ANRH2"
$ws.Cells.Item(25, 3).Value = "NA"
$ws.Cells.Item(25, 4).Value = "Nuclear Fuel"
$ws.Cells.Item(25, 5).Value = "h2"
$ws.Cells.Item(25, 6).Value = 80
$ws.Cells.Item(25, 7).Value = 10
$ws.Cells.Item(25, 11).Value = "Yes"
$ws.Cells.Item(25, 24).Value = "ANRH2"

# Row 26: iMSRHTSE
$ws.Cells.Item(26, 1).Value = "iMSRHTSE"
$ws.Cells.Item(26, 2).Value = "ANRH2"
$ws.Cells.Item(26, 3).Value = "NA"
$ws.Cells.Item(26, 4).Value = "Nuclear Fuel"
$ws.Cells.Item(26, 5).Value = "h2"
$ws.Cells.Item(26, 6).Value = 141
$ws.Cells.Item(26, 7).Value = 1
$ws.Cells.Item(26, 11).Value = "Yes"
$ws.Cells.Item(26, 24).Value = "ANRH2"

# Row 27: MicroHTSE
$ws.Cells.Item(27, 1).Value = "MicroHTSE"
$ws.Cells.Item(27, 2).Value = "ANRH2"
$ws.Cells.Item(27, 3).Value = "NA"
$ws.Cells.Item(27, 4).Value = "Nuclear Fuel"
$ws.Cells.Item(27, 5).Value = "h2"
$ws.Cells.Item(27, 6).Value = 6.7
$ws.Cells.Item(27, 7).Value = 50
$ws.Cells.Item(27, 11).Value = "Yes"
$ws.Cells.Item(27, 24).Value = "ANRH2"

# Match the author's final selection state on the sheet.
$ws.Range("B23:B27").Select()
